# Switched to the Malaysian site
# Update the "mostFrom" values (now mostly "Selangor", a Malaysian state, instead of
# "Indonesia" / "Mainland China") together with the refreshed avgPrice / avgSold
# figures that came from re-scraping the Malaysian Shopee site.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: shoes-for-men
$ws.Range("B2").Value = 33.6269767441861
$ws.Range("C2").Value = 32650.7209302326
$ws.Range("D2").Value = "Selangor"

# Row 3: shoes-for-women
$ws.Range("B3").Value = 19.8862790697674
$ws.Range("C3").Value = 46159.1860465116
$ws.Range("D3").Value = "Selangor"

# Row 4: shirts-for-men
$ws.Range("B4").Value = 15.7218181818182
$ws.Range("C4").Value = 39806.4848484848
$ws.Range("D4").Value = "Selangor"

# Row 5: dress-for-women
$ws.Range("B5").Value = 22.6553846153846
$ws.Range("C5").Value = 17571.1923076923
$ws.Range("D5").Value = "Indonesia"

# Row 6: bag-for-women
$ws.Range("B6").Value = 16.8791666666667
$ws.Range("C6").Value = 31939.1458333333
$ws.Range("D6").Value = "Selangor"
